# Insert a new data row at row 336 (pushing the existing row 336..397 down
# to 337..398) and populate it with the new Ají price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(336).Insert()

$ws.Range("A336").Value = 9
$ws.Range("B336").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C336").Value = "Metropolitana"
$ws.Range("D336").Value = 45034
$ws.Range("E336").Value = 13
$ws.Range("F336").Value = 100112021
$ws.Range("G336").Value = "Ají"
$ws.Range("H336").Value = "Americana (o)"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 34
$ws.Range("K336").Value = 12000
$ws.Range("L336").Value = 13000
$ws.Range("M336").Value = 12500
$ws.Range("N336").Value = "`$/caja 25 kilos"
$ws.Range("O336").Value = "Provincia de Limarí"
$ws.Range("P336").Value = 500
$ws.Range("Q336").Value = 25
$ws.Range("R336").Value = "Hortaliza"
